$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings like "26.997.35" that are not valid numbers
# (multiple "." separators) alongside ones that look numeric (e.g. "1.003").
# Force the whole Price column to Text so every new value is stored verbatim,
# matching the original authoring (inline strings), rather than letting Excel
# auto-coerce the numeric-looking ones into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @{
  2  = @{ D = "26.997.35"; E = "  +0.11%  " }
  3  = @{ D = "1.828.18";  E = "  +0.52%  " }
  4  = @{ D = "1.003";     E = "  -0.92%  " }
  5  = @{ D = "311.45";    E = "  +0.34%  " }
  6  = @{ D = "1.003";     E = "  -0.85%  " }
  7  = @{ D = "0.4629";    E = "  +0.06%  " }
  8  = @{ D = "0.3710";    E = "  +2.42%  " }
  9  = @{ D = "0.07341";   E = "  +0.73%  " }
  10 = @{ D = "0.8757";    E = "  +1.31%  " }
  11 = @{ D = "0.07876";   E = "  +3.76%  " }
  12 = @{ D = "19.75";     E = "  +0.12%  " }
  13 = @{ D = "1.849.56";  E = "  +1.16%  " }
  14 = @{ D = "5.337";     E = "  +0.35%  " }
  15 = @{ D = "6.555";     E = "  +1.82%  " }
  16 = @{ D = "91.29";     E = "  -1.68%  " }
  17 = @{ D = "1.006";     E = "  -0.57%  " }
  18 = @{ D = "0.000008829"; E = "  +2.41%  " }
  19 = @{ D = "1.003";     E = "  -0.71%  " }
  20 = @{ D = "14.79";     E = "  +2.59%  " }
  21 = @{ D = "27.005.98"; E = "  -0.84%  " }
  22 = @{ D = "5.096";     E = "  -1.07%  " }
  23 = @{ D = "10.53";     E = "  -0.19%  " }
  24 = @{ D = "2.069.70";  E = "  -1.36%  " }
  25 = @{ D = "152.75";    E = "  +1.07%  " }
  26 = @{ D = "1.857";     E = "  -0.22%  " }
  27 = @{ D = "18.46";     E = "  +1.45%  " }
  28 = @{ D = "2.043";     E = "  -2.39%  " }
  29 = @{ D = "5.125";     E = "  +1.33%  " }
  30 = @{ D = "115.57";    E = "  +0.26%  " }
  31 = @{ D = "0.08878";   E = "  -0.05%  " }
  32 = @{ D = "2.953";     E = "  -0.13%  " }
  33 = @{ D = "0.7281";    E = "  +0.44%  " }
  34 = @{ D = "4.435";     E = "  +0.68%  " }
  35 = @{ D = "1.131";     E = "  -0.08%  " }
  36 = @{ D = "2.468";     E = "  -1.65%  " }
  37 = @{ D = "0.01948";   E = "  +1.83%  " }
  38 = @{ D = "1.068";     E = "  -0.39%  " }
  39 = @{ D = "0.05219";   E = "  -0.47%  " }
  40 = @{ D = "2.947";     E = "  +0.77%  " }
  41 = @{ D = "7.112";     E = "  +0.16%  " }
  42 = @{ D = "0.5159";    E = "  -0.56%  " }
  43 = @{ D = "0.1622";    E = "  -0.24%  " }
  44 = @{ D = "8.172";     E = "  -0.32%  " }
  45 = @{ D = "0.4837";    E = "  -0.15%  " }
  46 = @{ D = "1.002";     E = "  -0.91%  " }
  47 = @{ D = "10.16";     E = "  +0.50%  " }
  48 = @{ D = "102.44";    E = "  -0.63%  " }
  49 = @{ D = "1.628";     E = "  -0.27%  " }
  50 = @{ D = "0.06198";   E = "  -0.62%  " }
  51 = @{ D = "64.94";     E = "  +1.10%  " }
}

foreach ($row in $updates.Keys) {
  $vals = $updates[$row]
  $ws.Range("D$row").Value = $vals.D
  $ws.Range("E$row").Value = $vals.E
}
